$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 9) duplicating the previous day's gold-price
# entry (row 8), reusing the same shared-string text and cell formatting
# (thin border on A9, thin border + wrap text on B9) already used by the
# other data rows.
$ws.Range("A9").Value = "26-09-2025"
$ws.Range("A9").Borders.LineStyle = 1

$rupee = [char]8377
$ws.Range("B9").Value = "The price of gold in India today is " + $rupee + "11,488 per gram for 24 karat gold, " + $rupee + "10,530 per gram for 22 karat gold and " + $rupee + "8,616 per gram for 18 karat gold (also called 999 gold)."
$ws.Range("B9").Borders.LineStyle = 1
$ws.Range("B9").WrapText = $true
